$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gphb5"
$ws.Range("C2").Value = "Tshr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3203896666666666
$ws.Range("H2").Value = 0.9611689999999999
$ws.Range("I2").Value = 0.07382438063517588
$ws.Range("J2").Value = 0.07382438063517589
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8749903333333333
$ws.Range("N2").Value = 2.624971
$ws.Range("O2").Value = 0.2670516933349977
$ws.Range("P2").Value = 0.2670516933349977
$ws.Range("Q2").Value = 0.2803378612332222
$ws.Range("R2").Value = 2.523040751099
$ws.Range("S2").Value = 0.01971492585803113
$ws.Range("T2").Value = 0.01971492585803114

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gphb5"
$ws.Range("C3").Value = "Tshr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3203896666666666
$ws.Range("H3").Value = 0.9611689999999999
$ws.Range("I3").Value = 0.07382438063517588
$ws.Range("J3").Value = 0.07382438063517589
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8147036666666666
$ws.Range("N3").Value = 2.444111
$ws.Range("O3").Value = 0.2486518827250642
$ws.Range("P3").Value = 0.2486518827250642
$ws.Range("Q3").Value = 0.2610226361954444
$ws.Range("R3").Value = 2.349203725759
$ws.Range("S3").Value = 0.01835657123594826
$ws.Range("T3").Value = 0.01835657123594826

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gphb5"
$ws.Range("C4").Value = "Tshr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3203896666666666
$ws.Range("H4").Value = 0.9611689999999999
$ws.Range("I4").Value = 0.07382438063517588
$ws.Range("J4").Value = 0.07382438063517589
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.586789
$ws.Range("N4").Value = 4.760367
$ws.Range("O4").Value = 0.484296423939938
$ws.Range("P4").Value = 0.484296423939938
$ws.Range("Q4").Value = 0.5083907987803332
$ws.Range("R4").Value = 4.575517189023
$ws.Range("S4").Value = 0.03575288354119649
$ws.Range("T4").Value = 0.0357528835411965

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gphb5"
$ws.Range("C5").Value = "Tshr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.195221666666667
$ws.Range("H5").Value = 9.585665
$ws.Range("I5").Value = 0.7362449076086343
$ws.Range("J5").Value = 0.7362449076086343
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8749903333333333
$ws.Range("N5").Value = 2.624971
$ws.Range("O5").Value = 0.2670516933349977
$ws.Range("P5").Value = 0.2670516933349977
$ws.Range("Q5").Value = 2.795788071190556
$ws.Range("R5").Value = 25.162092640715
$ws.Range("S5").Value = 0.1966154492861547
$ws.Range("T5").Value = 0.1966154492861547

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gphb5"
$ws.Range("C6").Value = "Tshr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.195221666666667
$ws.Range("H6").Value = 9.585665
$ws.Range("I6").Value = 0.7362449076086343
$ws.Range("J6").Value = 0.7362449076086343
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8147036666666666
$ws.Range("N6").Value = 2.444111
$ws.Range("O6").Value = 0.2486518827250642
$ws.Range("P6").Value = 0.2486518827250642
$ws.Range("Q6").Value = 2.603158807646111
$ws.Range("R6").Value = 23.428429268815
$ws.Range("S6").Value = 0.1830686824236279
$ws.Range("T6").Value = 0.1830686824236279

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gphb5"
$ws.Range("C7").Value = "Tshr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.195221666666667
$ws.Range("H7").Value = 9.585665
$ws.Range("I7").Value = 0.7362449076086343
$ws.Range("J7").Value = 0.7362449076086343
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.586789
$ws.Range("N7").Value = 4.760367
$ws.Range("O7").Value = 0.484296423939938
$ws.Range("P7").Value = 0.484296423939938
$ws.Range("Q7").Value = 5.070142593228333
$ws.Range("R7").Value = 45.631283339055
$ws.Range("S7").Value = 0.3565607758988516
$ws.Range("T7").Value = 0.3565607758988516

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gphb5"
$ws.Range("C8").Value = "Tshr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8242783333333333
$ws.Range("H8").Value = 2.472835
$ws.Range("I8").Value = 0.1899307117561898
$ws.Range("J8").Value = 0.1899307117561898
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8749903333333333
$ws.Range("N8").Value = 2.624971
$ws.Range("O8").Value = 0.2670516933349977
$ws.Range("P8").Value = 0.2670516933349977
$ws.Range("Q8").Value = 0.7212355736427778
$ws.Range("R8").Value = 6.491120162784999
$ws.Range("S8").Value = 0.05072131819081184
$ws.Range("T8").Value = 0.05072131819081185

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gphb5"
$ws.Range("C9").Value = "Tshr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8242783333333333
$ws.Range("H9").Value = 2.472835
$ws.Range("I9").Value = 0.1899307117561898
$ws.Range("J9").Value = 0.1899307117561898
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8147036666666666
$ws.Range("N9").Value = 2.444111
$ws.Range("O9").Value = 0.2486518827250642
$ws.Range("P9").Value = 0.2486518827250642
$ws.Range("Q9").Value = 0.6715425805205555
$ws.Range("R9").Value = 6.043883224685
$ws.Range("S9").Value = 0.04722662906548808
$ws.Range("T9").Value = 0.04722662906548809

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Gphb5"
$ws.Range("C10").Value = "Tshr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8242783333333333
$ws.Range("H10").Value = 2.472835
$ws.Range("I10").Value = 0.1899307117561898
$ws.Range("J10").Value = 0.1899307117561898
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.586789
$ws.Range("N10").Value = 4.760367
$ws.Range("O10").Value = 0.484296423939938
$ws.Range("P10").Value = 0.484296423939938
$ws.Range("Q10").Value = 1.307955792271666
$ws.Range("R10").Value = 11.771602130445
$ws.Range("S10").Value = 0.09198276449988985
$ws.Range("T10").Value = 0.09198276449988986
